try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.Worksheets.Item(4)   # 存款 (deposit) sheet

    # ------------------------------------------------------------------
    # 1. Turn row 1 (which used to be a stray duplicate data row) into a
    #    proper header row, and extend it with the new metadata columns
    #    G..M that every other sheet in the workbook already has.
    # ------------------------------------------------------------------
    $headerValues = @(
        "bank",               # B1
        "deposit_type",       # C1
        "currency",           # D1
        "owner",              # E1
        "total",              # F1
        "property_category",  # G1
        "category",           # H1
        "date",                # I1
        "legislator_name",    # J1
        "legislator_id",      # K1
        "source_file",        # L1
        "index"               # M1
    )

    # Make sure the whole header row (B1:M1) shares the same bold/boxed
    # style that B1:F1 already had, then fill in the values.
    $ws.Range("F1").Copy()
    $ws.Range("G1:M1").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    for ($i = 0; $i -lt $headerValues.Length; $i++) {
        $col = 2 + $i   # column B = 2
        $ws.Cells.Item(1, $col).Value = $headerValues[$i]
    }

    # ------------------------------------------------------------------
    # 2. Add the extra metadata columns (G..M) to the 10 existing data
    #    rows (rows 2-11). Columns A-F already hold the correct values.
    # ------------------------------------------------------------------
    $legislatorName = "許忠信"
    $legislatorId = 1749
    $sourceFile = "tmpa22c1"
    $category = "normal"
    $date = "2012-04-23"
    $propertyCategory = "deposit"

    $ws.Range("F2").Copy()
    $ws.Range("G2:M11").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    for ($r = 2; $r -le 11; $r++) {
        $indexValue = $ws.Cells.Item($r, 1).Value()

        $ws.Cells.Item($r, 7).Value = $propertyCategory    # G: property_category
        $ws.Cells.Item($r, 8).Value = $category             # H: category
        $ws.Cells.Item($r, 9).Value = $date                 # I: date
        $ws.Cells.Item($r, 10).Value = $legislatorName      # J: legislator_name
        $ws.Cells.Item($r, 11).Value = $legislatorId        # K: legislator_id
        $ws.Cells.Item($r, 12).Value = $sourceFile          # L: source_file
        $ws.Cells.Item($r, 13).Value = $indexValue          # M: index
    }

    Write-Host "deposit sheet updated"
} catch {
    Write-Host "ERROR:" $_.Exception.Message
}
